$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Guard columns D and E (which are stored as inline/shared TEXT strings, even
# though many of the new values look numeric) against Excel's automatic
# "string that looks like a number becomes a number" coercion on
# Range.Value assignment. We temporarily force a Text number format so the
# assignment is stored as a literal string, then clear the formatting again
# (restoring the default, un-styled cell) once all values are written so the
# cells end up with their original (style-less) appearance.
$guard = $ws.Range("D2:E51")
$guard.NumberFormat = "@"

$ws.Range("D2").Value = "62.438.99"
$ws.Range("E2").Value = "  -2.73%  "
$ws.Range("D3").Value = "3.049.71"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "531.28"
$ws.Range("E5").Value = "  -5.45%  "
$ws.Range("D6").Value = "131.92"
$ws.Range("E6").Value = "  -6.33%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.041.98"
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").Value = "6.18"
$ws.Range("E11").Value = "  -8.44%  "
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").Value = "33.94"
$ws.Range("E14").Value = "  -6.00%  "
$ws.Range("D15").Value = "3.517.99"
$ws.Range("E15").Value = "  -3.73%  "
$ws.Range("D16").Value = "62.290.33"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "0.110"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Value = "3.021.31"
$ws.Range("E18").Value = "  -4.13%  "
$ws.Range("D19").Value = "6.54"
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("D20").Value = "477.55"
$ws.Range("E20").Value = "  -6.44%  "
$ws.Range("D21").Value = "13.14"
$ws.Range("E21").Value = "  -5.47%  "
$ws.Range("D22").Value = "0.688"
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("D23").Value = "6.99"
$ws.Range("E23").Value = "  -5.06%  "
$ws.Range("D24").Value = "78.50"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "11.96"
$ws.Range("E25").Value = "  -5.55%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "2.65"
$ws.Range("E27").Value = "  -6.14%  "
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  -6.35%  "
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "25.65"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("D31").Value = "1.84"
$ws.Range("E31").Value = "  -11.95%  "
$ws.Range("D32").Value = "1.10"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "2.35"
$ws.Range("E33").Value = "  -8.74%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "56.48"
$ws.Range("E34").Value = "  +4.58%  "
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").Value = "470.07"
$ws.Range("E37").Value = "  -15.50%  "
$ws.Range("D38").Value = "3.075.51"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "0.0389"
$ws.Range("E39").Value = "  -8.35%  "
$ws.Range("D40").Value = "0.0783"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("D41").Value = "0.114"
$ws.Range("E41").Value = "  -5.46%  "
$ws.Range("D42").Value = "7.99"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("D43").Value = "2.57"
$ws.Range("E43").Value = "  -5.17%  "
$ws.Range("D45").Value = "0.247"
$ws.Range("E45").Value = "  -5.42%  "
$ws.Range("D46").Value = "0.0₃0547"
$ws.Range("E46").Value = "  +6.80%  "
$ws.Range("D47").Value = "2.02"
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("D48").Value = "119.31"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "24.19"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.107"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").Value = "2.32"
$ws.Range("E51").Value = "  +8.89%  "

$guard.ClearFormats()
